$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.921940299144911
$ws.Range("C2").Value = 5.087245312770772
$ws.Range("D2").Value = 10.73728626441013
$ws.Range("F2").Value = 33.25604469951789
$ws.Range("G2").Value = 3.656053689898157
$ws.Range("I2").Value = 22.47910669656748
$ws.Range("J2").Value = 11.3796198932624
$ws.Range("K2").Value = 10.03814764600269
$ws.Range("O2").Value = 24.49909318008374

$ws.Range("B3").Value = 9.602431641741099
$ws.Range("C3").Value = 4.855981155922448
$ws.Range("D3").Value = 10.6268897115429
$ws.Range("F3").Value = 33.32047813878802
$ws.Range("G3").Value = 3.65788762994707
$ws.Range("I3").Value = 22.58825463905698
$ws.Range("J3").Value = 11.36116819666873
$ws.Range("K3").Value = 9.819035691043004
$ws.Range("O3").Value = 24.59612689651433

$ws.Range("B4").Value = 9.401864274649773
$ws.Range("C4").Value = 4.708529636081979
$ws.Range("D4").Value = 10.56064387513844
$ws.Range("F4").Value = 33.36846784632593
$ws.Range("G4").Value = 3.659073094977194
$ws.Range("I4").Value = 22.6596610653076
$ws.Range("J4").Value = 11.35210217758756
$ws.Range("K4").Value = 9.683172433075352
$ws.Range("O4").Value = 24.66098613734018

$ws.Range("B5").Value = 9.319151430931703
$ws.Range("C5").Value = 4.647147431292371
$ws.Range("D5").Value = 10.53406270417903
$ws.Range("F5").Value = 33.39013786892627
$ws.Range("G5").Value = 3.659571169788117
$ws.Range("I5").Value = 22.68986320510084
$ws.Range("J5").Value = 11.34897951253314
$ws.Range("K5").Value = 9.627547540010763
$ws.Range("O5").Value = 24.68874197657875

$ws.Range("B6").Value = 9.305361636748387
$ws.Range("C6").Value = 4.636879251323616
$ws.Range("D6").Value = 10.5296747546231
$ws.Range("F6").Value = 33.39386366186236
$ws.Range("G6").Value = 3.659654781361944
$ws.Range("I6").Value = 22.69494488634605
$ws.Range("J6").Value = 11.34849560149896
$ws.Range("K6").Value = 9.618297619711171
$ws.Range("O6").Value = 24.69343077430645

$ws.Range("B7").Value = 9.400752578011854
$ws.Range("C7").Value = 4.707706945562177
$ws.Range("D7").Value = 10.5602836775808
$ws.Range("F7").Value = 33.36875154454219
$ws.Range("G7").Value = 3.659079751439052
$ws.Range("I7").Value = 22.66006391444017
$ws.Range("J7").Value = 11.35205774580283
$ws.Range("K7").Value = 9.682423206884092
$ws.Range("O7").Value = 24.66135510015043

$ws.Range("B8").Value = 9.812755261027704
$ws.Range("C8").Value = 5.008680309741789
$ws.Range("D8").Value = 10.69891919776174
$ws.Range("F8").Value = 33.27650927595106
$ws.Range("G8").Value = 3.656673728441713
$ws.Range("I8").Value = 22.51582978525059
$ws.Range("J8").Value = 11.37278970579855
$ws.Range("K8").Value = 9.962917121402992
$ws.Range("O8").Value = 24.53145297913323

$ws.Range("B9").Value = 10.58097209770023
$ws.Range("C9").Value = 5.552762261653081
$ws.Range("D9").Value = 10.98173346374754
$ws.Range("F9").Value = 33.16269856210133
$ws.Range("G9").Value = 3.652424838237589
$ws.Range("I9").Value = 22.26782407325422
$ws.Range("J9").Value = 11.43126888985285
$ws.Range("K9").Value = 10.49931529633936
$ws.Range("O9").Value = 24.31872452096871

$ws.Range("B10").Value = 11.11531225801192
$ws.Range("C10").Value = 5.921226049061376
$ws.Range("D10").Value = 11.19448011913829
$ws.Range("F10").Value = 33.12021789586265
$ws.Range("G10").Value = 3.64958627796544
$ws.Range("I10").Value = 22.1068672342915
$ws.Range("J10").Value = 11.48489198106684
$ws.Range("K10").Value = 10.88115339242874
$ws.Range("O10").Value = 24.18819674654653

$ws.Range("B11").Value = 11.35082705644589
$ws.Range("C11").Value = 6.081555159517225
$ws.Range("D11").Value = 11.29197294879144
$ws.Range("F11").Value = 33.10986128912306
$ws.Range("G11").Value = 3.648355775967352
$ws.Range("I11").Value = 22.03826317384504
$ws.Range("J11").Value = 11.5115476656727
$ws.Range("K11").Value = 11.05144362735292
$ws.Range("O11").Value = 24.13444316307361

$ws.Range("B12").Value = 11.4388489726549
$ws.Range("C12").Value = 6.141186611945034
$ws.Range("D12").Value = 11.32896265944442
$ws.Range("F12").Value = 33.10723068679175
$ws.Range("G12").Value = 3.647898507614061
$ws.Range("I12").Value = 22.01294900801853
$ws.Range("J12").Value = 11.521961576216
$ws.Range("K12").Value = 11.11538009815986
$ws.Range("O12").Value = 24.11489951654572

$ws.Range("B13").Value = 11.41994468755143
$ws.Range("C13").Value = 6.128392499839538
$ws.Range("D13").Value = 11.32099358746249
$ws.Range("F13").Value = 33.10773979005147
$ws.Range("G13").Value = 3.647996602458696
$ws.Range("I13").Value = 22.01837128829685
$ws.Range("J13").Value = 11.5197046122687
$ws.Range("K13").Value = 11.1016355266673
$ws.Range("O13").Value = 24.11907244525721

$ws.Range("B14").Value = 11.35809237084986
$ws.Range("C14").Value = 6.086482927707404
$ws.Range("D14").Value = 11.29501486384005
$ws.Range("F14").Value = 33.10961898368498
$ws.Range("G14").Value = 3.648317982188117
$ws.Range("I14").Value = 22.03616723502279
$ws.Range("J14").Value = 11.51239804346083
$ws.Range("K14").Value = 11.05671503676618
$ws.Range("O14").Value = 24.1328190118796

$ws.Range("B15").Value = 11.32005256656855
$ws.Range("C15").Value = 6.060670327537558
$ws.Range("D15").Value = 11.27911050652988
$ws.Range("F15").Value = 33.11093823245867
$ws.Range("G15").Value = 3.648515967761779
$ws.Range("I15").Value = 22.0471543587411
$ws.Range("J15").Value = 11.50796406482988
$ws.Range("K15").Value = 11.02912681413018
$ws.Range("O15").Value = 24.14134497007723

$ws.Range("B16").Value = 11.0997620567712
$ws.Range("C16").Value = 5.910598507296209
$ws.Range("D16").Value = 11.18812040677716
$ws.Range("F16").Value = 33.12107534320569
$ws.Range("G16").Value = 3.649667913388462
$ws.Range("I16").Value = 22.11144363467798
$ws.Range("J16").Value = 11.48319504591618
$ws.Range("K16").Value = 10.86995082209426
$ws.Range("O16").Value = 24.19182309283141

$ws.Range("B17").Value = 10.96262844159705
$ws.Range("C17").Value = 5.816642001788022
$ws.Range("D17").Value = 11.13246128584086
$ws.Range("F17").Value = 33.1295924953883
$ws.Range("G17").Value = 3.650390129460152
$ws.Range("I17").Value = 22.15206593369452
$ws.Range("J17").Value = 11.46857581725572
$ws.Range("K17").Value = 10.77138528553297
$ws.Range("O17").Value = 24.22423236028949

$ws.Range("B18").Value = 10.88304524205522
$ws.Range("C18").Value = 5.761916451788431
$ws.Range("D18").Value = 11.10051663161964
$ws.Range("F18").Value = 33.13533540918835
$ws.Range("G18").Value = 3.650811251867383
$ws.Range("I18").Value = 22.1758652272758
$ws.Range("J18").Value = 11.46038050525115
$ws.Range("K18").Value = 10.71437527662225
$ws.Range("O18").Value = 24.24340251369036

$ws.Range("B19").Value = 10.8559806522697
$ws.Range("C19").Value = 5.743270926741745
$ws.Range("D19").Value = 11.0897135288604
$ws.Range("F19").Value = 33.13742476801106
$ws.Range("G19").Value = 3.650954820897033
$ws.Range("I19").Value = 22.183997844328
$ws.Range("J19").Value = 11.45764250032822
$ws.Range("K19").Value = 10.69501994662889
$ws.Range("O19").Value = 24.24998399912827

$ws.Range("B20").Value = 10.97730035363913
$ws.Range("C20").Value = 5.826714913251678
$ws.Range("D20").Value = 11.13837938052753
$ws.Range("F20").Value = 33.12859845707711
$ws.Range("G20").Value = 3.650312656348198
$ws.Range("I20").Value = 22.1476966530823
$ws.Range("J20").Value = 11.47011002254786
$ws.Range("K20").Value = 10.78191107425474
$ws.Range("O20").Value = 24.22072754690186

$ws.Range("B21").Value = 11.37629200220121
$ws.Range("C21").Value = 6.098822383757665
$ws.Range("D21").Value = 11.30264375763841
$ws.Range("F21").Value = 33.10903196696338
$ws.Range("G21").Value = 3.64822334950408
$ws.Range("I21").Value = 22.03092208654897
$ws.Range("J21").Value = 11.51453552016509
$ws.Range("K21").Value = 11.06992461124027
$ws.Range("O21").Value = 24.12875926051202

$ws.Range("B22").Value = 11.63025151810597
$ws.Range("C22").Value = 6.270343492272305
$ws.Range("D22").Value = 11.41040130718862
$ws.Range("F22").Value = 33.10377049257245
$ws.Range("G22").Value = 3.646908534179331
$ws.Range("I22").Value = 21.95847802886257
$ws.Range("J22").Value = 11.54543248593466
$ws.Range("K22").Value = 11.25493543734408
$ws.Range("O22").Value = 24.07338501969999

$ws.Range("B23").Value = 11.49535390108403
$ws.Range("C23").Value = 6.179387196310312
$ws.Range("D23").Value = 11.35286260531708
$ws.Range("F23").Value = 33.10588966563376
$ws.Range("G23").Value = 3.647605653930951
$ws.Range("I23").Value = 21.99678791252742
$ws.Range("J23").Value = 11.52877365709071
$ws.Range("K23").Value = 11.15650477585587
$ws.Range("O23").Value = 24.10250533012982

$ws.Range("B24").Value = 10.97066949146864
$ws.Range("C24").Value = 5.822163154812436
$ws.Range("D24").Value = 11.13570363772732
$ws.Range("F24").Value = 33.12904522576554
$ws.Range("G24").Value = 3.650347663527504
$ws.Range("I24").Value = 22.14967061876224
$ws.Range("J24").Value = 11.46941575544031
$ws.Range("K24").Value = 10.77715343020242
$ws.Range("O24").Value = 24.22231039852848

$ws.Range("B25").Value = 10.37804452492168
$ws.Range("C25").Value = 5.4108705192409
$ws.Range("D25").Value = 10.90422799099638
$ws.Range("F25").Value = 33.18627812164793
$ws.Range("G25").Value = 3.653524344447297
$ws.Range("I25").Value = 22.33118528935575
$ws.Range("J25").Value = 11.41356020884882
$ws.Range("K25").Value = 10.35608700582952
$ws.Range("O25").Value = 24.37175982719633
